$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2/D2 to numeric values
$ws.Range("C2").Value = 149948
$ws.Range("D2").Value = 1

# Update C3/D3 to numeric values
$ws.Range("C3").Value = 237310
$ws.Range("D3").Value = 2

# Delete column F (month_year) entirely, shifting dimension back to A1:E3
$ws.Range("F1:F3").Delete()
